$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (C) column for every existing data row (2..390)
#    from 2023-09-12 (45181) to 2023-09-13 (45182).
$ws.Range("C2:C390").Value = 45182

# 2) Row 390 gains an explicit custom row height (15pt) in the saved XML.
$ws.Rows.Item(390).RowHeight = 15

# 3) Append a new data row (391) for case "A 42432-2023".
$ws.Range("A391").Value = "A 42432-2023"

$ws.Range("B391").Value = 45180
$ws.Range("B391").NumberFormat = "YYYY-MM-DD"

$ws.Range("C391").Value = 45182
$ws.Range("C391").NumberFormat = "YYYY-MM-DD"

$ws.Range("D391").Value = "SKÅNE LÄN"
$ws.Range("E391").Value = "ÖSTRA GÖINGE"

$ws.Range("G391").Value = 0.9
$ws.Range("H391").Value = 0
$ws.Range("I391").Value = 0
$ws.Range("J391").Value = 0
$ws.Range("K391").Value = 0
$ws.Range("L391").Value = 0
$ws.Range("M391").Value = 0
$ws.Range("N391").Value = 0
$ws.Range("O391").Value = 0
$ws.Range("P391").Value = 0
$ws.Range("Q391").Value = 0

# R391 stays empty but keeps the wrap-text styling used throughout column R.
$ws.Range("R391").WrapText = $true
